$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that already carries the timestamp number format used by
# column G (s="2" -> numFmtId 165, "YYYY-MM-DD HH:MM:SS").
$tsFormat = $ws.Cells.Item(66, 7).NumberFormat

$rows = @(
    @{ A = 111; B = "PT301224170951"; C = "debolina"; D = 25; E = "mci";  F = "Normal";   G = 45656.71918981482; H = 29; I = "May be Normal" },
    @{ A = 114; B = "PT301224171900"; C = "q";        D = 4;  E = $null;  F = $null;      G = 45656.72152777778; H = 0;  I = "Severe" },
    @{ A = 115; B = "PT301224171947"; C = "f";        D = 5;  E = $null;  F = $null;      G = 45656.72207175926; H = 0;  I = "Severe" },
    @{ A = 116; B = "PT301224172258"; C = "e";        D = 2;  E = $null;  F = $null;      G = 45656.72428240741; H = 0;  I = "Severe" }
)

$startRow = 67
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D

    if ($null -ne $row.E) {
        $ws.Cells.Item($r, 5).Value = $row.E
    }
    if ($null -ne $row.F) {
        $ws.Cells.Item($r, 6).Value = $row.F
    }

    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 7).NumberFormat = $tsFormat

    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
